# Natmi LR-pairs result for Thbs2-Cd47, regenerated following Dr Hou's advice:
# the analysis now also considers "ECs" as a possible sending cluster, so the
# table grows from 2 sending clusters x 3 target clusters (6 rows) to
# 3 sending clusters x 3 target clusters (9 rows). Header row (row 1) is
# unchanged; rows 2-10 are (re)written with the refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Thbs2"
$ws.Cells.Item(2,3).Value = "Cd47"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.8911683333333333
$ws.Cells.Item(2,8).Value = 2.673505
$ws.Cells.Item(2,9).Value = 0.02693425114262819
$ws.Cells.Item(2,10).Value = 0.02693425114262819
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 61.27353666666666
$ws.Cells.Item(2,14).Value = 183.82061
$ws.Cells.Item(2,15).Value = 0.3474604587406809
$ws.Cells.Item(2,16).Value = 0.3474604587406808
$ws.Cells.Item(2,17).Value = 54.60503554867222
$ws.Cells.Item(2,18).Value = 491.44531993805
$ws.Cells.Item(2,19).Value = 0.009358587257854298
$ws.Cells.Item(2,20).Value = 0.009358587257854297

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Thbs2"
$ws.Cells.Item(3,3).Value = "Cd47"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.8911683333333333
$ws.Cells.Item(3,8).Value = 2.673505
$ws.Cells.Item(3,9).Value = 0.02693425114262819
$ws.Cells.Item(3,10).Value = 0.02693425114262819
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 89.66709900000001
$ws.Cells.Item(3,14).Value = 269.001297
$ws.Cells.Item(3,15).Value = 0.5084702637939138
$ws.Cells.Item(3,16).Value = 0.5084702637939138
$ws.Cells.Item(3,17).Value = 79.908479170665
$ws.Cells.Item(3,18).Value = 719.1763125359851
$ws.Cells.Item(3,19).Value = 0.01369526578358368
$ws.Cells.Item(3,20).Value = 0.01369526578358368

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Thbs2"
$ws.Cells.Item(4,3).Value = "Cd47"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.8911683333333333
$ws.Cells.Item(4,8).Value = 2.673505
$ws.Cells.Item(4,9).Value = 0.02693425114262819
$ws.Cells.Item(4,10).Value = 0.02693425114262819
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 25.40615466666667
$ws.Cells.Item(4,14).Value = 76.218464
$ws.Cells.Item(4,15).Value = 0.1440692774654054
$ws.Cells.Item(4,16).Value = 0.1440692774654053
$ws.Cells.Item(4,17).Value = 22.64116051070222
$ws.Cells.Item(4,18).Value = 203.77044459632
$ws.Cells.Item(4,19).Value = 0.003880398101190212
$ws.Cells.Item(4,20).Value = 0.003880398101190211

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Thbs2"
$ws.Cells.Item(5,3).Value = "Cd47"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 24.359699
$ws.Cells.Item(5,8).Value = 73.07909699999999
$ws.Cells.Item(5,9).Value = 0.7362360466408275
$ws.Cells.Item(5,10).Value = 0.7362360466408276
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 61.27353666666666
$ws.Cells.Item(5,14).Value = 183.82061
$ws.Cells.Item(5,15).Value = 0.3474604587406809
$ws.Cells.Item(5,16).Value = 0.3474604587406808
$ws.Cells.Item(5,17).Value = 1492.604909865463
$ws.Cells.Item(5,18).Value = 13433.44418878917
$ws.Cells.Item(5,19).Value = 0.2558129145072472
$ws.Cells.Item(5,20).Value = 0.2558129145072472

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Thbs2"
$ws.Cells.Item(6,3).Value = "Cd47"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 24.359699
$ws.Cells.Item(6,8).Value = 73.07909699999999
$ws.Cells.Item(6,9).Value = 0.7362360466408275
$ws.Cells.Item(6,10).Value = 0.7362360466408276
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 89.66709900000001
$ws.Cells.Item(6,14).Value = 269.001297
$ws.Cells.Item(6,15).Value = 0.5084702637939138
$ws.Cells.Item(6,16).Value = 0.5084702637939138
$ws.Cells.Item(6,17).Value = 2184.263541843201
$ws.Cells.Item(6,18).Value = 19658.37187658881
$ws.Cells.Item(6,19).Value = 0.3743541368500498
$ws.Cells.Item(6,20).Value = 0.3743541368500499

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Thbs2"
$ws.Cells.Item(7,3).Value = "Cd47"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 24.359699
$ws.Cells.Item(7,8).Value = 73.07909699999999
$ws.Cells.Item(7,9).Value = 0.7362360466408275
$ws.Cells.Item(7,10).Value = 0.7362360466408276
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 25.40615466666667
$ws.Cells.Item(7,14).Value = 76.218464
$ws.Cells.Item(7,15).Value = 0.1440692774654054
$ws.Cells.Item(7,16).Value = 0.1440692774654053
$ws.Cells.Item(7,17).Value = 618.8862804274452
$ws.Cells.Item(7,18).Value = 5569.976523847007
$ws.Cells.Item(7,19).Value = 0.1060689952835305
$ws.Cells.Item(7,20).Value = 0.1060689952835305

$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Thbs2"
$ws.Cells.Item(8,3).Value = "Cd47"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 7.835938333333334
$ws.Cells.Item(8,8).Value = 23.507815
$ws.Cells.Item(8,9).Value = 0.2368297022165442
$ws.Cells.Item(8,10).Value = 0.2368297022165442
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 61.27353666666666
$ws.Cells.Item(8,14).Value = 183.82061
$ws.Cells.Item(8,15).Value = 0.3474604587406809
$ws.Cells.Item(8,16).Value = 0.3474604587406808
$ws.Cells.Item(8,17).Value = 480.1356547852389
$ws.Cells.Item(8,18).Value = 4321.22089306715
$ws.Cells.Item(8,19).Value = 0.08228895697557931
$ws.Cells.Item(8,20).Value = 0.0822889569755793

$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Thbs2"
$ws.Cells.Item(9,3).Value = "Cd47"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 7.835938333333334
$ws.Cells.Item(9,8).Value = 23.507815
$ws.Cells.Item(9,9).Value = 0.2368297022165442
$ws.Cells.Item(9,10).Value = 0.2368297022165442
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 89.66709900000001
$ws.Cells.Item(9,14).Value = 269.001297
$ws.Cells.Item(9,15).Value = 0.5084702637939138
$ws.Cells.Item(9,16).Value = 0.5084702637939138
$ws.Cells.Item(9,17).Value = 702.6258582928951
$ws.Cells.Item(9,18).Value = 6323.632724636056
$ws.Cells.Item(9,19).Value = 0.1204208611602803
$ws.Cells.Item(9,20).Value = 0.1204208611602803

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Thbs2"
$ws.Cells.Item(10,3).Value = "Cd47"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 7.835938333333334
$ws.Cells.Item(10,8).Value = 23.507815
$ws.Cells.Item(10,9).Value = 0.2368297022165442
$ws.Cells.Item(10,10).Value = 0.2368297022165442
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 25.40615466666667
$ws.Cells.Item(10,14).Value = 76.218464
$ws.Cells.Item(10,15).Value = 0.1440692774654054
$ws.Cells.Item(10,16).Value = 0.1440692774654053
$ws.Cells.Item(10,17).Value = 199.0810612551289
$ws.Cells.Item(10,18).Value = 1791.72955129616
$ws.Cells.Item(10,19).Value = 0.03411988408068464
$ws.Cells.Item(10,20).Value = 0.03411988408068464

